$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.808.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.42%  '

$ws.Range("D3").Value = "'2.355.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.17%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'240.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '

$ws.Range("D6").Value = "'0.670"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.34%  '

$ws.Range("D7").Value = "'73.66"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.76%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("D9").Value = "'0.603"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.45%  '

$ws.Range("E10").Value = '  +1.64%  '

$ws.Range("D11").Value = "'60.88"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.29%  '

$ws.Range("D12").Value = "'35.59"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +10.45%  '

$ws.Range("E13").Value = '  +0.52%  '

$ws.Range("D14").Value = "'7.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.09%  '

$ws.Range("D15").Value = "'16.28"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.09%  '

$ws.Range("D16").Value = "'0.918"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.13%  '

$ws.Range("D17").Value = "'2.351.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").Value = "'43.764.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.28%  '

$ws.Range("E19").Value = '  +1.06%  '

$ws.Range("D20").Value = "'77.85"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.12%  '

$ws.Range("D21").Value = "'6.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.66%  '

$ws.Range("D22").Value = "'253.59"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.80%  '

$ws.Range("E23").Value = '  +0.01%  '

$ws.Range("E24").Value = '  +3.13%  '

$ws.Range("D25").Value = "'1.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.40%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").Value = "'10.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.66%  '

$ws.Range("E28").Value = '  +0.73%  '

$ws.Range("D29").Value = "'175.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.12%  '

$ws.Range("D30").Value = "'22.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.04%  '

$ws.Range("E31").Value = '  +0.46%  '

$ws.Range("D32").Value = "'0.134"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.82%  '

$ws.Range("E33").Value = '  -2.78%  '

$ws.Range("D34").Value = "'5.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.44%  '

$ws.Range("D35").Value = "'5.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.34%  '

$ws.Range("D36").Value = "'3.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.65%  '

$ws.Range("D37").Value = "'6.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.72%  '

$ws.Range("E38").Value = '  +2.11%  '

$ws.Range("D39").Value = "'0.0278"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.89%  '

$ws.Range("D40").Value = "'5.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +13.68%  '

$ws.Range("D41").Value = "'65.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +11.65%  '

$ws.Range("D42").Value = "'20.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.25%  '

$ws.Range("B43").Value = 'Cronos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D43").Value = "'0.107"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.78%  '

$ws.Range("D44").Value = "'9.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.36%  '

$ws.Range("B45").Value = 'Algorand'
$ws.Range("C45").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D45").Value = "'0.203"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.61%  '

$ws.Range("D46").Value = "'2.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.33%  '

$ws.Range("E47").Value = '  +0.07%  '

$ws.Range("D48").Value = "'1.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.72%  '

$ws.Range("E49").Value = '  -1.95%  '

$ws.Range("D50").Value = "'98.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.20%  '

$ws.Range("E51").Value = '  +2.06%  '
